# Update "想去人数" (wish-to-attend count) figures in the "展览" and
# "全部类型" sheets to reflect the newly generated output data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 537
$wsExhibit.Range("F3").Value = 6409
$wsExhibit.Range("F4").Value = 397
$wsExhibit.Range("F6").Value = 135
$wsExhibit.Range("F8").Value = 74
$wsExhibit.Range("F9").Value = 571
$wsExhibit.Range("F10").Value = 44

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 537
$wsAll.Range("F3").Value = 6409
$wsAll.Range("F4").Value = 397
$wsAll.Range("F7").Value = 135
$wsAll.Range("F10").Value = 74
$wsAll.Range("F11").Value = 571
$wsAll.Range("F12").Value = 44
